$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 43 edits -----------------------------------------------------

# C43: "Em atendimento" -> "Aguardando atendimento externo" (style s="6" kept as-is)
$ws.Range("C43").Value = "Aguardando atendimento externo"

# L43: "Registrado" -> "Respondido" (style s="6" kept as-is)
$ws.Range("L43").Value = "Respondido"

# T43 was an empty, unstyled cell; it needs to become a wrapped-text cell
# (style s="1", same format as S43/U43) with new content. Clone the format
# from U43 (already s="1") and then set the value.
$ws.Range("U43").Copy()
$ws.Range("T43").PasteSpecial(-4122)
$ws.Range("T43").Value = "19/07/2021 11:26 Kaique Ferreira Henrique de Souza: `nME - Planejamento e Desenvolvimento - Sicoob Uni`nPrezados(as), boa tarde..."

# U43: replace text content, keep existing style s="1"
$ws.Range("U43").Value = "19/07/2021 11:26 Kaique Ferreira Henrique de Souza: `n<img src='/icons/progresstrail.mail/16' style='vertical-align: text-bott..."

# --- New row 83 ---------------------------------------------------------
# Clone per-column formatting from the existing last data row (82) so the
# new row reuses the same style indices instead of creating new ones.

$ws.Range("A82:F82").Copy()
$ws.Range("A83").PasteSpecial(-4122)

$ws.Range("I82:L82").Copy()
$ws.Range("I83").PasteSpecial(-4122)

$ws.Range("M82").Copy()
$ws.Range("M83").PasteSpecial(-4122)

$ws.Range("Q82").Copy()
$ws.Range("Q83").PasteSpecial(-4122)

$ws.Range("R82").Copy()
$ws.Range("R83").PasteSpecial(-4122)

$ws.Range("S82").Copy()
$ws.Range("S83").PasteSpecial(-4122)

$ws.Range("U82").Copy()
$ws.Range("U83").PasteSpecial(-4122)

$ws.Range("A83").Value = "I2107-176424"
$ws.Range("B83").Value = "Marcelo Da Silva"
$ws.Range("C83").Value = "Aguardando atendimento"
$ws.Range("D83").Value = "5145 - PAC:00 - MEDCRED RIBEIRÃO PRETO"
$ws.Range("E83").Value = "2015 - Desenvolvimento Estratégico"
$ws.Range("F83").Value = "2015 - PAC: 00 - SICOOB UNI"
$ws.Range("G83").Value = $false
$ws.Range("H83").Value = $false
$ws.Range("I83").Value = "2015. Inteligência de Negócios"
$ws.Range("J83").Value = "Sobras"
$ws.Range("K83").Value = "2015 - Desenvolvimento Estratégico"
$ws.Range("L83").Value = "Registrado"
$ws.Range("M83").Value = 44396.741006944445
$ws.Range("Q83").Value = 0.0
$ws.Range("R83").Value = "Central - 2015 - Central Sicoob Uni"
$ws.Range("S83").Value = "19/07/2021 17:47 Atualizado pelo Sistema: `nCooperativa: 5145 - PAC:00 - MEDCRED DE RIBEIRÃO PRETO`nTelefone: 16 36256900`n19/0..."
$ws.Range("U83").Value = "19/07/2021 17:48 Portal de Servi&ccedil;os do CCS &lt;no-reply@sicoob.com.br&gt;: `n<img src='/icons/progresstrail.mail/16' st..."
